$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H17").Value = 2292.8728
$ws.Range("J17").Value = 2292.8728
$ws.Range("L17").Value = 6878.6184
$ws.Range("N17").Value = -7214.6184

$ws.Range("H19").Value = 319.7143
$ws.Range("I19").Value = 300.33334
$ws.Range("J19").Value = 325
$ws.Range("K19").Value = 300.33334
$ws.Range("L19").Value = 325
$ws.Range("M19").Value = -125.33334
$ws.Range("N19").Value = -675

$ws.Range("H32").Value = 27778994
$ws.Range("I32").Value = 166666900
$ws.Range("J32").Value = 1414.3
$ws.Range("K32").Value = 166666900
$ws.Range("L32").Value = 1414.3
$ws.Range("M32").Value = -166666574
$ws.Range("N32").Value = -2066.3

$ws.Range("H33").Value = 248.41096
$ws.Range("I33").Value = 197.67606
$ws.Range("K33").Value = 197.67606
$ws.Range("M33").Value = 31.32393999999999

$ws.Range("H62").Value = 4237.696
$ws.Range("I62").Value = 6135.4546
$ws.Range("J62").Value = 2498.0833
$ws.Range("K62").Value = 6135.4546
$ws.Range("L62").Value = 2498.0833
$ws.Range("M62").Value = -5511.4546
$ws.Range("N62").Value = -3746.0833

$ws.Range("H65").Value = 4237.696
$ws.Range("I65").Value = 6135.4546
$ws.Range("J65").Value = 2498.0833
$ws.Range("K65").Value = 30677.273
$ws.Range("L65").Value = 12490.4165
$ws.Range("M65").Value = -27557.273
$ws.Range("N65").Value = -18730.4165

$ws.Range("H76").Value = 2992.0232
$ws.Range("I76").Value = 2882.12
$ws.Range("J76").Value = 3144.6667
$ws.Range("K76").Value = 2882.12
$ws.Range("L76").Value = 3144.6667
$ws.Range("M76").Value = -2567.12
$ws.Range("N76").Value = -3774.6667

$ws.Range("H79").Value = 2992.0232
$ws.Range("I79").Value = 2882.12
$ws.Range("J79").Value = 3144.6667
$ws.Range("K79").Value = 2882.12
$ws.Range("L79").Value = 3144.6667
$ws.Range("M79").Value = -1790.12
$ws.Range("N79").Value = -5328.6667

$ws.Range("H129").Value = 1055.3776
$ws.Range("I129").Value = 2329
$ws.Range("J129").Value = 972.31525
$ws.Range("K129").Value = 6987
$ws.Range("L129").Value = 2916.94575
$ws.Range("M129").Value = -1987
$ws.Range("N129").Value = -12916.94575

$ws.Range("H137").Value = 2082948.4
$ws.Range("I137").Value = 6993919.5
$ws.Range("J137").Value = 5229.8076
$ws.Range("K137").Value = 20981758.5
$ws.Range("L137").Value = 15689.4228
$ws.Range("M137").Value = -20979208.5
$ws.Range("N137").Value = -20789.4228

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 12220.38
$ws.Range("I32").Value = 11227.387
$ws.Range("J32").Value = 19502.334
$ws.Range("K32").Value = 11227.387
$ws.Range("L32").Value = 19502.334
$ws.Range("M32").Value = -10940.387
$ws.Range("N32").Value = -20076.334

$ws.Range("H122").Value = 1738.5
$ws.Range("I122").Value = 1598.125
$ws.Range("J122").Value = 2300
$ws.Range("K122").Value = 4794.375
$ws.Range("L122").Value = 6900
$ws.Range("M122").Value = -2344.375
$ws.Range("N122").Value = -11800

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H134").Value = 2596.0625
$ws.Range("I134").Value = 1635.6177
$ws.Range("J134").Value = 3684.5667
$ws.Range("K134").Value = 4906.8531
$ws.Range("L134").Value = 11053.7001
$ws.Range("M134").Value = -2371.8531
$ws.Range("N134").Value = -16123.7001

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H16").Value = 1414.5625
$ws.Range("I16").Value = 1317.9231
$ws.Range("J16").Value = 1833.3334
$ws.Range("K16").Value = 1317.9231
$ws.Range("L16").Value = 1833.3334
$ws.Range("M16").Value = -1030.9231
$ws.Range("N16").Value = -2407.3334

$ws.Range("H31").Value = 4427
$ws.Range("I31").Value = 2265.85
$ws.Range("K31").Value = 2265.85
$ws.Range("M31").Value = -1970.85

$ws.Range("H34").Value = 4427
$ws.Range("I34").Value = 2265.85
$ws.Range("K34").Value = 2265.85
$ws.Range("M34").Value = -2063.85

$ws.Range("H58").Value = 2095.5
$ws.Range("I58").Value = 1426.7778
$ws.Range("J58").Value = 3098.5833
$ws.Range("K58").Value = 1426.7778
$ws.Range("L58").Value = 3098.5833
$ws.Range("M58").Value = -1223.7778
$ws.Range("N58").Value = -3504.5833

$ws.Range("H86").Value = 4237.7144
$ws.Range("J86").Value = 4499.6665
$ws.Range("L86").Value = 4499.6665
$ws.Range("N86").Value = -6745.6665

$ws.Range("H89").Value = 4237.7144
$ws.Range("J89").Value = 4499.6665
$ws.Range("L89").Value = 22498.3325
$ws.Range("N89").Value = -33730.3325

$ws.Range("H99").Value = 3417.3044
$ws.Range("I99").Value = 3104
$ws.Range("J99").Value = 3759.0908
$ws.Range("K99").Value = 3104
$ws.Range("L99").Value = 3759.0908
$ws.Range("M99").Value = -1606
$ws.Range("N99").Value = -6755.0908

$ws.Range("H113").Value = 1414.5625
$ws.Range("I113").Value = 1317.9231
$ws.Range("J113").Value = 1833.3334
$ws.Range("K113").Value = 1317.9231
$ws.Range("L113").Value = 1833.3334
$ws.Range("M113").Value = 852.0769
$ws.Range("N113").Value = -6173.3334

$ws.Range("H126").Value = 3417.3044
$ws.Range("I126").Value = 3104
$ws.Range("J126").Value = 3759.0908
$ws.Range("K126").Value = 9312
$ws.Range("L126").Value = 11277.2724
$ws.Range("M126").Value = -6842
$ws.Range("N126").Value = -16217.2724

$ws.Range("H134").Value = 1033767.7
$ws.Range("I134").Value = 1566013.1
$ws.Range("J134").Value = 235399.5
$ws.Range("K134").Value = 4698039.300000001
$ws.Range("L134").Value = 706198.5
$ws.Range("M134").Value = -4695504.300000001
$ws.Range("N134").Value = -711268.5

$ws.Range("H136").Value = 2095.5
$ws.Range("I136").Value = 1426.7778
$ws.Range("J136").Value = 3098.5833
$ws.Range("K136").Value = 4280.3334
$ws.Range("L136").Value = 9295.749899999999
$ws.Range("M136").Value = -1730.3334
$ws.Range("N136").Value = -14395.7499

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H113").Value = 772.4
$ws.Range("J113").Value = 904.9091
$ws.Range("L113").Value = 2714.7273
$ws.Range("N113").Value = -7054.7273

$ws.Range("H132").Value = 3089.75
$ws.Range("I132").Value = 1099.6
$ws.Range("K132").Value = 9896.4
$ws.Range("M132").Value = -7366.4

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H22").Value = 750
$ws.Range("I22").Value = 550
$ws.Range("J22").Value = 1150
$ws.Range("K22").Value = 550
$ws.Range("L22").Value = 1150
$ws.Range("M22").Value = -255
$ws.Range("N22").Value = -1740

$ws.Range("H27").Value = 750
$ws.Range("I27").Value = 550
$ws.Range("J27").Value = 1150
$ws.Range("K27").Value = 550
$ws.Range("L27").Value = 1150
$ws.Range("M27").Value = -443
$ws.Range("N27").Value = -1364

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H126").Value = 1716.174
$ws.Range("I126").Value = 1191.1177
$ws.Range("J126").Value = 3203.8333
$ws.Range("K126").Value = 3573.3531
$ws.Range("L126").Value = 9611.499899999999
$ws.Range("M126").Value = -1103.3531
$ws.Range("N126").Value = -14551.4999

$ws.Range("H132").Value = 1404333.8
$ws.Range("I132").Value = 1891658.8
$ws.Range("J132").Value = 3274.5
$ws.Range("K132").Value = 5674976.4
$ws.Range("L132").Value = 9823.5
$ws.Range("M132").Value = -5672446.4
$ws.Range("N132").Value = -14883.5

Write-Host "Applied scheduled runner updates to ALC, ARM, BSM, CRP, CUL, LTW, WVR sheets"
